# Apply the "fixed segment and powertrain split" edit.
#
# The core of this change extends the linear extrapolation of the
# Aluminium-content-per-powertrain-segment trend (columns DR:EW, i.e. the
# years after the last "real" data point in DQ) on the "Data" sheet for
# each of the 5 data rows (2-6). Each new cell continues the straight-line
# trend of the two preceding cells: new = prev + prev - prevprev.
#
# It also reproduces the reachable cosmetic changes: the workbook's active
# sheet/tab moves from "Cover" to "Data" (with an updated selection), and
# every sheet gets an explicit A4/portrait page setup.

$wb = $excel.ActiveWorkbook

$wsCover     = $wb.Worksheets.Item("Cover")
$wsData      = $wb.Worksheets.Item("Data")
$wsDataOld   = $wb.Worksheets.Item("DataOld")
$wsAlContent = $wb.Worksheets.Item("Al_Content")
$wsInflow    = $wb.Worksheets.Item("Inflow")

# --- Extend the linear trend forward from DQ into DR:EW for rows 2-6 -----
# DR{row} = DQ{row} + DQ{row} - DP{row}
# DS{row}:EW{row} = previous cell + previous cell - cell before that
foreach ($row in 2..6) {
    $wsData.Range("DR$row").FormulaR1C1 = "=RC[-1]+RC[-1]-RC[-2]"
    $wsData.Range("DS$row`:EW$row").FormulaR1C1 = "=RC[-1]+RC[-1]-RC[-2]"
}

# --- Recalculate so the cached <v> values land in the saved workbook -----
$excel.CalculateFull()

# --- Give every sheet an explicit page setup (A4 / portrait) -------------
foreach ($ws in @($wsCover, $wsData, $wsDataOld, $wsAlContent, $wsInflow)) {
    $ws.PageSetup.PaperSize = 9
    $ws.PageSetup.Orientation = 1
}

# --- Switch the active/selected tab from Cover to Data, updating the -----
# --- view's selected cell to match the new area of interest --------------
$wsData.Activate()
$wsData.Range("EW7").Select() | Out-Null
